# "Updated Digits Kinematics Length Data"
# The DP_Length (mm) column (column A) values for rows 2-6 were recomputed;
# update them to the new measurements.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 29.360877984147542
$ws.Range("A3").Value = 21.110617731369206
$ws.Range("A4").Value = 20.401408548431156
$ws.Range("A5").Value = 20.154365035892347
$ws.Range("A6").Value = 18.460373398173722
